# N2200_database - column title changing.xlsx
#
# This edit fills in the newly-introduced "gate_material" column (AI) for
# every existing data row with the electrode metal that was actually used
# (Au / Al / Cu), un-highlights the AI1 header cell (it was marked with a
# yellow "new column" fill that is no longer needed now that it is
# populated), normalizes a stray cell-style left over on BK36:BN45, and
# moves the on-screen selection over to the newly edited area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fill column AI ("gate_material") for rows 2-157.
#    The values were entered in contiguous blocks; rows 65-68 and 151
#    belong to entries that are out of scope for this pass and are left
#    untouched (they still have no gate_material value).
# ---------------------------------------------------------------------
$gateBlocks = @(
    @{ Start = 2;   End = 14;  Metal = "Au" },
    @{ Start = 15;  End = 64;  Metal = "Al" },
    @{ Start = 69;  End = 99;  Metal = "Al" },
    @{ Start = 100; End = 101; Metal = "Au" },
    @{ Start = 102; End = 109; Metal = "Cu" },
    @{ Start = 110; End = 135; Metal = "Al" },
    @{ Start = 136; End = 150; Metal = "Au" },
    @{ Start = 152; End = 157; Metal = "Al" }
)

foreach ($block in $gateBlocks) {
    $rng = $ws.Range("AI$($block.Start):AI$($block.End)")
    $rng.Value = $block.Metal
}

# ---------------------------------------------------------------------
# 2) AI1 ("gate_material" header) no longer needs the yellow "new
#    column" highlight now that the column has been populated - clear
#    its fill while keeping the sheet's standard vertical-center
#    alignment. Do this by lifting the (already unhighlighted) format
#    from A1, which carries exactly that combination.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AI1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) BK36:BN45 were left with a leftover "applyFill" style from an
#    earlier pass; normalize them to match the rest of the column
#    (same look as row 46 and below) by copying that format over.
# ---------------------------------------------------------------------
$ws.Range("BK46:BN46").Copy() | Out-Null
$ws.Range("BK36:BN45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Move the visible window / selection to the area that was just
#    edited.
# ---------------------------------------------------------------------
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("AC1") } catch {}
try { $excel.ActiveWindow.Zoom = 70 } catch {}
$ws.Range("AI10").Select() | Out-Null
